$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. The last bullet ("text-align, text-decoration 속성 ") loses its
#    trailing space.
# ---------------------------------------------------------------------
$null = $d.Content.Find.Execute(
    "text-align, text-decoration 속성 ", $true, $false, $false, $false, $false,
    $true, 1, $false, "text-align, text-decoration 속성", 2)

# ---------------------------------------------------------------------
# 2. Append the new material that was typed after that bullet:
#       (blank line)
#       (blank line)
#       2022-02-21
#       10강. css 속성 3
#       position 속성        <- new bullet, same list as the bullet above
# ---------------------------------------------------------------------
$wNs = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

$blank      = "<w:p $wNs/>"
$dateLine   = "<w:p $wNs>" +
                "<w:r><w:rPr><w:rFonts w:hint=`"eastAsia`"/></w:rPr><w:lastRenderedPageBreak/><w:t>2</w:t></w:r>" +
                "<w:r><w:t>022-02-21</w:t></w:r>" +
              "</w:p>"
$titleLine  = "<w:p $wNs>" +
                "<w:r><w:rPr><w:rFonts w:hint=`"eastAsia`"/></w:rPr><w:t>1</w:t></w:r>" +
                "<w:r><w:t>0</w:t></w:r>" +
                "<w:r><w:rPr><w:rFonts w:hint=`"eastAsia`"/></w:rPr><w:t>강.</w:t></w:r>" +
                "<w:r><w:t xml:space=`"preserve`"> </w:t></w:r>" +
                "<w:proofErr w:type=`"spellStart`"/>" +
                "<w:r><w:rPr><w:rFonts w:hint=`"eastAsia`"/></w:rPr><w:t>c</w:t></w:r>" +
                "<w:r><w:t>ss</w:t></w:r>" +
                "<w:proofErr w:type=`"spellEnd`"/>" +
                "<w:r><w:t xml:space=`"preserve`"> </w:t></w:r>" +
                "<w:r><w:rPr><w:rFonts w:hint=`"eastAsia`"/></w:rPr><w:t>속성 3</w:t></w:r>" +
              "</w:p>"
$bulletLine = "<w:p $wNs>" +
                "<w:pPr>" +
                  "<w:pStyle w:val=`"a4`"/>" +
                  "<w:numPr><w:ilvl w:val=`"0`"/><w:numId w:val=`"2`"/></w:numPr>" +
                  "<w:ind w:leftChars=`"0`"/>" +
                  "<w:rPr><w:rFonts w:hint=`"eastAsia`"/></w:rPr>" +
                "</w:pPr>" +
                "<w:r><w:rPr><w:rFonts w:hint=`"eastAsia`"/></w:rPr><w:t>p</w:t></w:r>" +
                "<w:r><w:t>osition</w:t></w:r>" +
                "<w:r><w:t xml:space=`"preserve`"> 속성</w:t></w:r>" +
              "</w:p>"

$xml = $blank + $blank + $dateLine + $titleLine + $bulletLine

# Insert strictly after the very last paragraph mark of the document body
# (before sectPr), so the existing content is left untouched.
$endRange = $d.Range($d.Content.End, $d.Content.End)
$null = $endRange.InsertXML($xml)

Write-Output "Done. Paragraphs: $($d.Paragraphs.Count)"
